$d = $word.ActiveDocument

# 1. Rename "Alissa Vanderbelt" -> "Alissa McQuate"
$d.Content.Find.Execute("Alissa Vanderbelt", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alissa McQuate", 2)

# 2. Insert a new paragraph "Alexis Whitacre" immediately before "Logan Whitacre"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Logan Whitacre`r") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newPara.Range.Text = "Alexis Whitacre"
        break
    }
}

# 3. Remove the old "Alexis Zenich" paragraph (it moved later in the collection
#    after the insertion above, so re-scan rather than relying on old indices).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Alexis Zenich`r") {
        $p.Range.Delete()
        break
    }
}
